$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 4 by shifting rows 4..8 down to 5..9 ---
# Work bottom-up so we never clobber a row before we've read it.
for ($r = 8; $r -ge 4; $r--) {
    $dest = $r + 1
    $v = $ws.Range("A$r").Value2
    $h = $ws.Rows.Item($r).RowHeight
    $ws.Range("A$dest").Value = $v
    $ws.Rows.Item($dest).RowHeight = $h
}

# --- New content for B3 (registered first so it lands at the lower shared-string index) ---
$ws.Range("B3").Value = "2022/3/8完成"

# --- New content for row 4 (the freshly inserted row) ---
$ws.Range("A4").Value = "第一层嵌套规定不能用GridLayoutGroup，解除之前的逻辑"
$ws.Rows.Item(4).RowHeight = 46.2

# --- Apply word-wrap to every used cell so they all land on one shared style ---
$ws.Range("A1").WrapText = $true
$ws.Range("B1").WrapText = $true
$ws.Range("A2").WrapText = $true
$ws.Range("B2").WrapText = $true
$ws.Range("A3").WrapText = $true
$ws.Range("B3").WrapText = $true
$ws.Range("A4").WrapText = $true
$ws.Range("A5").WrapText = $true
$ws.Range("A6").WrapText = $true
$ws.Range("A7").WrapText = $true
$ws.Range("A8").WrapText = $true
$ws.Range("A9").WrapText = $true

# --- Update the selected cell shown in the saved view ---
$ws.Range("J3").Select() | Out-Null
